$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1543.2667
$ws.Range("I28").Value = 1024.5454
$ws.Range("K28").Value = 1024.5454
$ws.Range("M28").Value = -539.5454

$ws.Range("H33").Value = 329.41666
$ws.Range("I33").Value = 329.41666
$ws.Range("K33").Value = 329.41666
$ws.Range("M33").Value = -100.41666

$ws.Range("H51").Value = 9199.6
$ws.Range("I51").Value = 9199.6
$ws.Range("K51").Value = 9199.6
$ws.Range("M51").Value = -8715.6

$ws.Range("H64").Value = 4999.1665

$ws.Range("H67").Value = 4999.1665

$ws.Range("H96").Value = 1532
$ws.Range("I96").Value = 1772
$ws.Range("K96").Value = 5316
$ws.Range("M96").Value = -3943

$ws.Range("H98").Value = 754.5
$ws.Range("I98").Value = 754.5
$ws.Range("K98").Value = 754.5
$ws.Range("M98").Value = 743.5

$ws.Range("H116").Value = 5999.75
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558

$ws.Range("H122").Value = 754.5
$ws.Range("I122").Value = 754.5
$ws.Range("K122").Value = 2263.5
$ws.Range("M122").Value = 186.5

$ws.Range("H135").Value = 1011.4286
$ws.Range("I135").Value = 1011.4286
$ws.Range("K135").Value = 9102.857399999999
$ws.Range("M135").Value = -6567.857399999999

$ws.Range("H137").Value = 1845.1177
$ws.Range("I137").Value = 1898
$ws.Range("J137").Value = 999
$ws.Range("K137").Value = 5694
$ws.Range("L137").Value = 2997
$ws.Range("M137").Value = -3144
$ws.Range("N137").Value = -8097

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 2
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 171

$ws.Range("H32").Value = 4232.1
$ws.Range("I32").Value = 3076.4443
$ws.Range("K32").Value = 3076.4443
$ws.Range("M32").Value = -2789.4443

$ws.Range("H61").Value = 2142.4285
$ws.Range("J61").Value = 1998
$ws.Range("L61").Value = 1998
$ws.Range("N61").Value = -2422

$ws.Range("H132").Value = 1856.375
$ws.Range("I132").Value = 1856.375
$ws.Range("K132").Value = 5569.125
$ws.Range("M132").Value = -3039.125

$ws.Range("H136").Value = 2142.4285
$ws.Range("J136").Value = 1998
$ws.Range("L136").Value = 5994
$ws.Range("N136").Value = -11094

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2448.4348
$ws.Range("I86").Value = 2426.4285
$ws.Range("K86").Value = 2426.4285
$ws.Range("M86").Value = -1303.4285

$ws.Range("H89").Value = 2448.4348
$ws.Range("I89").Value = 2426.4285
$ws.Range("K89").Value = 12132.1425
$ws.Range("M89").Value = -6516.1425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1130.3077
$ws.Range("J107").Value = 1369.2222
$ws.Range("L107").Value = 1369.2222
$ws.Range("N107").Value = -5209.2222

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2000
$ws.Range("I51").Value = 2000
$ws.Range("K51").Value = 6000
$ws.Range("M51").Value = -5540

$ws.Range("H55").Value = 2621.4285
$ws.Range("I55").Value = 150
$ws.Range("K55").Value = 450
$ws.Range("M55").Value = -273

$ws.Range("H80").Value = 3333
$ws.Range("I80").Value = 3333
$ws.Range("K80").Value = 9999
$ws.Range("M80").Value = -9063

$ws.Range("H83").Value = 3333
$ws.Range("I83").Value = 3333
$ws.Range("K83").Value = 29997
$ws.Range("M83").Value = -25317

$ws.Range("H92").Value = 280.8
$ws.Range("I92").Value = 235.66667
$ws.Range("K92").Value = 707.00001
$ws.Range("M92").Value = 540.99999

$ws.Range("H109").Value = 1922.1666
$ws.Range("J109").Value = 2561.4285
$ws.Range("L109").Value = 7684.2855
$ws.Range("N109").Value = -9764.2855

$ws.Range("H121").Value = 909813.6
$ws.Range("J121").Value = 1111850
$ws.Range("L121").Value = 3335550
$ws.Range("N121").Value = -3338170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3995
$ws.Range("I41").Value = 3995
$ws.Range("K41").Value = 3995
$ws.Range("M41").Value = -3640

$ws.Range("H70").Value = 5656.4
$ws.Range("I70").Value = 5966.4287
$ws.Range("J70").Value = 4933
$ws.Range("K70").Value = 5966.4287
$ws.Range("L70").Value = 4933
$ws.Range("M70").Value = -5696.4287
$ws.Range("N70").Value = -5473

$ws.Range("H73").Value = 5656.4
$ws.Range("I73").Value = 5966.4287
$ws.Range("J73").Value = 4933
$ws.Range("K73").Value = 5966.4287
$ws.Range("L73").Value = 4933
$ws.Range("M73").Value = -5030.4287
$ws.Range("N73").Value = -6805

$ws.Range("H80").Value = 2938.6
$ws.Range("I80").Value = 2771.5
$ws.Range("K80").Value = 2771.5
$ws.Range("M80").Value = -1773.5

$ws.Range("H83").Value = 2938.6
$ws.Range("I83").Value = 2771.5
$ws.Range("K83").Value = 13857.5
$ws.Range("M83").Value = -8865.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 522500
$ws.Range("I2").Value = 45000
$ws.Range("K2").Value = 45000
$ws.Range("M2").Value = -44888

$ws.Range("H4").Value = 5962.625
$ws.Range("I4").Value = 2002
$ws.Range("J4").Value = 6528.4287
$ws.Range("K4").Value = 2002
$ws.Range("L4").Value = 6528.4287
$ws.Range("M4").Value = -1889
$ws.Range("N4").Value = -6754.4287

$ws.Range("H41").Value = 19659.334
$ws.Range("J41").Value = 17489.5
$ws.Range("L41").Value = 17489.5
$ws.Range("N41").Value = -18269.5

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws.Range("H136").Value = 7671.7144
$ws.Range("J136").Value = 9666.666999999999
$ws.Range("L136").Value = 29000.001
$ws.Range("N136").Value = -34100.001
